# MAS Report templates formatting update
#
# Content-level changes captured by this script (the rest of the upstream
# diff - fileVersion/rupBuild, absPath, revisionPtr GUID, x14ac:dyDescent
# hints, xr:uid stamps - are Excel-build/save-environment metadata that is
# regenerated by whichever Excel instance performs the save and isn't
# something a COM script should hand-author):
#
#   1. Clear the sample/placeholder figures that shipped in the template
#      (B6:C7 - "1.11/1001" and "2.99/2000") so the form is blank again.
#   2. Give the two value cells (B6, B7) an explicit currency number
#      format ("$"#,##0.00) - this is the new numFmtId 164 added to
#      styles.xml.
#   3. Update the active/selected cell saved in the sheet view to B15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clear the placeholder sample values for the two data rows.
$ws.Range("B6:C7").ClearContents()

# 2. Apply a currency number format to the value cells.
$ws.Range("B6").NumberFormat = """$""#,##0.00"
$ws.Range("B7").NumberFormat = """$""#,##0.00"

# 3. Move/save the selection to B15.
$ws.Range("B15").Select()
